$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.455.07"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.981.25"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'381.48"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").Value = "'104.15"
$ws.Range("E6").Value = "  +2.99%  "
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D13").Value = "3.455.97"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "'18.45"
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").Value = "'7.80"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").Value = "2.983.22"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "'11.14"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "51.470.66"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "'70.23"
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("D24").Value = "'266.93"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("E25").Value = "  +2.58%  "
$ws.Range("E26").Value = "  -3.85%  "
$ws.Range("D27").Value = "'7.29"
$ws.Range("E27").Value = "  -3.64%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'26.04"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").Value = "'10.43"
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("D33").Value = "'34.66"
$ws.Range("E33").Value = "  +3.61%  "
$ws.Range("D34").Value = "'51.36"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").Value = "'2.07"
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'3.28"
$ws.Range("E38").Value = "  +3.38%  "
$ws.Range("D39").Value = "'16.94"
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("E40").Value = "  +4.18%  "
$ws.Range("D41").Value = "'0.116"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("D43").Value = "'3.81"
$ws.Range("E43").Value = "  +11.64%  "
$ws.Range("D44").Value = "'125.55"
$ws.Range("E44").Value = "  +4.07%  "
$ws.Range("D45").Value = "'21.40"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "'2.03"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("D49").Value = "2.022.11"
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("D50").Value = "3.279.11"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "'0.0330"
$ws.Range("E51").Value = "  +0.47%  "
